# Ajout de fonctionallites accepter et refuser don avec envoie d'un courriel
#
# On the "SprintBacklog2" sheet, the three tasks belonging to the user story
# "traitement de don par employe" (rows 18-20: the UI for the employee's
# donation list, the accept/refuse logic, and the donor e-mail notification)
# move from the "a faire" (to-do, column E) status to "en cours" (in
# progress, column F) on the little status tracker. The round yellow marker
# ("bullet") together with its centered formatting moves from column E to
# column F on each of those rows, and the now-empty E cell takes on the
# plain (non-centered) formatting that F used to have. The row heights also
# settle to new values as the content reflows, and the active selection on
# that sheet moves to J17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SprintBacklog2")

# Keep this sheet the active / selected tab, as it was before the edit.
$ws.Activate()

$scratch = $ws.Range("Z1")

foreach ($r in 18, 19, 20) {
    $eCell = $ws.Range("E$r")
    $fCell = $ws.Range("F$r")

    # Swap the cell *formats* between E and F (E currently has the
    # centered "active" look, F has the plain "inactive" look) via
    # copy / paste-special so the existing style entries get reused
    # instead of new ones being synthesized.
    $eCell.Copy()
    $scratch.PasteSpecial(-4122)   # xlPasteFormats

    $fCell.Copy()
    $eCell.PasteSpecial(-4122)

    $scratch.Copy()
    $fCell.PasteSpecial(-4122)

    $scratch.Clear()

    # Move the status marker value itself from column E ("a faire") to
    # column F ("en cours").
    $marker = $eCell.Value()
    $fCell.Value = $marker
    $eCell.Value = $null
}

# Row heights settle to their new auto values after the content move.
$ws.Rows.Item(18).RowHeight = 13.8
$ws.Rows.Item(19).RowHeight = 23.95
$ws.Rows.Item(20).RowHeight = 13.8

# Active cell / selection moves to J17.
$ws.Range("J17").Select()
